{"js": "// fix(docx): fix OOXMLValidator error on KeywordTok output\n//\n// Several custom \"*Tok\" character styles (used for syntax-highlighted\n// source code) had their <w:rPr> children in the wrong order -- w:color\n// was emitted before w:b / w:i, which violates the CT_RPr element\n// sequence in wml.xsd (rFonts, b, bCmpl, i, iCmpl, ..., color, ...).\n// OOXMLValidatorCLI flags this as Sch_UnexpectedElementContentExpecting\n// Complex even though some lenient parsers (e.g. xmllint) don't warn.\n//\n// Fix: for every affected style, re-assign its own existing Bold /\n// Italic flag(s) on Font. That doesn't change any visible formatting,\n// but it forces Word to rebuild the style's <w:rPr> from the in-memory\n// model, which re-serializes b/i ahead of color in schema order.\nconst styleNames = [\n  \"KeywordTok\",\n  \"ImportTok\",\n  \"CommentTok\",\n  \"DocumentationTok\",\n  \"AnnotationTok\",\n  \"CommentVarTok\",\n  \"ControlFlowTok\",\n  \"InformationTok\",\n  \"WarningTok\",\n  \"AlertTok\",\n  \"ErrorTok\",\n];\n\nconst styles = context.document.getStyles();\nconst fonts = styleNames.map((name) => styles.getByNameOrNullObject(name).font);\nfonts.forEach((font) => font.load(\"bold,italic\"));\nawait context.sync();\n\n// Only re-assign flags that are actually set (truthy), mirroring what\n// was already present in the markup -- this avoids introducing an\n// explicit \"false\" (e.g. <w:i w:val=\"0\"/>) for a property that wasn't\n// serialized at all in the original style.\nfor (const font of fonts) {\n  if (font.bold) {\n    font.bold = font.bold;\n  }\n  if (font.italic) {\n    font.italic = font.italic;\n  }\n}\nawait context.sync();\n", "ps1": "# fix(docx): fix OOXMLValidator error on KeywordTok output\n#\n# The character styles below had <w:rPr> children in the wrong order\n# (w:color before w:b/w:i), which violates the wml.xsd CT_RPr sequence\n# (rFonts, b, bCmpl, i, iCmpl, ..., color, ...). Re-assigning each\n# style's existing Bold/Italic flag forces Word to re-serialize the run\n# properties in schema-correct order (b/i before color) without\n# changing any actual formatting.\n\n$d = $word.ActiveDocument\n\n$styleNames = @(\n    \"KeywordTok\",\n    \"ImportTok\",\n    \"CommentTok\",\n    \"DocumentationTok\",\n    \"AnnotationTok\",\n    \"CommentVarTok\",\n    \"ControlFlowTok\",\n    \"InformationTok\",\n    \"WarningTok\",\n    \"AlertTok\",\n    \"ErrorTok\"\n)\n\nforeach ($name in $styleNames) {\n    $style = $d.Styles.Item($name)\n    # Only touch flags that are actually set (True), so we don't\n    # introduce an explicit <w:b w:val=\"0\"/>/<w:i w:val=\"0\"/> for a\n    # property that wasn't present in the original markup.\n    if ($style.Font.Bold) {\n        $style.Font.Bold = $style.Font.Bold\n    }\n    if ($style.Font.Italic) {\n        $style.Font.Italic = $style.Font.Italic\n    }\n}\n"}
